$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 91

$ws.Cells.Item($newRow, 1).Value = "'2025-10-21"
$ws.Cells.Item($newRow, 1).Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = "'15:22:34"
$ws.Cells.Item($newRow, 2).Style = "Normal"

$ws.Cells.Item($newRow, 3).Value = "1.00 EUR = 1,747.5624"
$ws.Cells.Item($newRow, 3).Style = "Normal"
